$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 0.8262640128117136
$ws.Range("C6").Value = 0.8361151140958082
$ws.Range("D6").Value = 0.8262640128117136
$ws.Range("E6").Value = 0.827942324055343
$ws.Range("F6").Value = 0.813406543125143
$ws.Range("G6").Value = 0.8226640132472383
$ws.Range("H6").Value = 0.813406543125143
$ws.Range("I6").Value = 0.8142274241189149
$ws.Range("J6").Value = 0.7618165179592772
$ws.Range("K6").Value = 0.7715510414696543
$ws.Range("L6").Value = 0.7618165179592772
$ws.Range("M6").Value = 0.7609529102436011
$ws.Range("N6").Value = 0.8068862960420956
$ws.Range("O6").Value = 0.8228227651648335
$ws.Range("P6").Value = 0.8068862960420956
$ws.Range("Q6").Value = 0.8082327911765731
$ws.Range("R6").Value = 0.8068862960420956
$ws.Range("S6").Value = 0.8211608850485754
$ws.Range("T6").Value = 0.8068862960420956
$ws.Range("U6").Value = 0.8079642493184324
$ws.Range("V6").Value = 0.8176161061541981
$ws.Range("W6").Value = 0.8206326696415698
$ws.Range("X6").Value = 0.8176161061541981
$ws.Range("Y6").Value = 0.8169122273468202
$ws.Range("B7").Value = 0.8520018302447954
$ws.Range("C7").Value = 0.8560606784491798
$ws.Range("D7").Value = 0.8520018302447954
$ws.Range("E7").Value = 0.8526172370702561
$ws.Range("F7").Value = 0.8712880347746511
$ws.Range("G7").Value = 0.8747590861838775
$ws.Range("H7").Value = 0.8712880347746511
$ws.Range("I7").Value = 0.8715642118514711
$ws.Range("J7").Value = 0.8283687943262411
$ws.Range("K7").Value = 0.8405181121986848
$ws.Range("L7").Value = 0.8283687943262411
$ws.Range("M7").Value = 0.8305318399463862
$ws.Range("N7").Value = 0.8797986730725235
$ws.Range("O7").Value = 0.8840739423079615
$ws.Range("P7").Value = 0.8797986730725235
$ws.Range("Q7").Value = 0.8795942503981061
$ws.Range("V7").Value = 0.8625714939373141
$ws.Range("W7").Value = 0.8692726364377155
$ws.Range("X7").Value = 0.8625714939373141
$ws.Range("Y7").Value = 0.86305094578302
